$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- 1) The empty paragraph right before the bookmark paragraph gains <w:ilvl w:val="0"/> ---
$n = $d.Paragraphs.Count
$blankPara = $d.Paragraphs.Item($n - 1)
$blankXml = "<w:p $wNs>" +
    "<w:pPr>" +
        "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='0'/></w:numPr>" +
        "<w:rPr><w:rFonts w:hint='default'/><w:lang w:val='pt-BR'/></w:rPr>" +
    "</w:pPr>" +
"</w:p>"
$blankPara.Range.InsertXML($blankXml) | Out-Null

# --- 2) The final paragraph (holding the _GoBack bookmark) becomes a "Final notes" heading ---
$n = $d.Paragraphs.Count
$bmPara = $d.Paragraphs.Item($n)
$headingXml = "<w:p $wNs>" +
    "<w:pPr>" +
        "<w:pStyle w:val='2'/>" +
        "<w:bidi w:val='0'/>" +
        "<w:rPr><w:rFonts w:hint='default'/><w:lang w:val='pt-BR'/></w:rPr>" +
    "</w:pPr>" +
    "<w:r>" +
        "<w:rPr><w:rFonts w:hint='default'/><w:lang w:val='pt-BR'/></w:rPr>" +
        "<w:t>Final notes</w:t>" +
    "</w:r>" +
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
    "<w:bookmarkEnd w:id='0'/>" +
"</w:p>"
$bmPara.Range.InsertXML($headingXml) | Out-Null

# --- 3) Append 5 new bulleted (numId=2) list paragraphs after the heading paragraph ---
for ($i = 0; $i -lt 5; $i++) {
    $lastPara = $d.Paragraphs.Last
    $lastPara.Range.InsertParagraphAfter() | Out-Null
}

# Helper to build the standard list-item pPr used by every new bullet paragraph below.
$listPPr = "<w:pPr>" +
        "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr>" +
        "<w:ind w:left='420' w:leftChars='0' w:hanging='420' w:firstLineChars='0'/>" +
        "<w:rPr><w:rFonts w:hint='default'/><w:lang w:val='pt-BR'/></w:rPr>" +
    "</w:pPr>"

$n = $d.Paragraphs.Count
$p1 = $d.Paragraphs.Item($n - 4)
$xml1 = "<w:p $wNs>" + $listPPr +
    "<w:r>" +
        "<w:rPr><w:rFonts w:hint='default'/><w:lang w:val='pt-BR'/></w:rPr>" +
        "<w:t>notice any overlapping subproblems</w:t>" +
    "</w:r>" +
"</w:p>"
$p1.Range.InsertXML($xml1) | Out-Null

$p2 = $d.Paragraphs.Item($n - 3)
$xml2 = "<w:p $wNs>" + $listPPr +
    "<w:r>" +
        "<w:rPr><w:rFonts w:hint='default'/><w:lang w:val='pt-BR'/></w:rPr>" +
        "<w:t>decide what is the triveally smallest input</w:t>" +
    "</w:r>" +
"</w:p>"
$p2.Range.InsertXML($xml2) | Out-Null

$p3 = $d.Paragraphs.Item($n - 2)
$xml3 = "<w:p $wNs>" + $listPPr +
    "<w:r>" +
        "<w:rPr><w:rFonts w:hint='default'/><w:lang w:val='pt-BR'/></w:rPr>" +
        "<w:t xml:space='preserve'>think recursively to use </w:t>" +
    "</w:r>" +
    "<w:r>" +
        "<w:rPr><w:rFonts w:hint='default'/><w:b/><w:bCs/><w:lang w:val='pt-BR'/></w:rPr>" +
        "<w:t>memoization</w:t>" +
    "</w:r>" +
"</w:p>"
$p3.Range.InsertXML($xml3) | Out-Null

$p4 = $d.Paragraphs.Item($n - 1)
$xml4 = "<w:p $wNs>" + $listPPr +
    "<w:r>" +
        "<w:rPr><w:rFonts w:hint='default'/><w:lang w:val='pt-BR'/></w:rPr>" +
        "<w:t xml:space='preserve'>think iteratively to use </w:t>" +
    "</w:r>" +
    "<w:r>" +
        "<w:rPr><w:rFonts w:hint='default'/><w:b/><w:bCs/><w:lang w:val='pt-BR'/></w:rPr>" +
        "<w:t>tabulation</w:t>" +
    "</w:r>" +
"</w:p>"
$p4.Range.InsertXML($xml4) | Out-Null

$p5 = $d.Paragraphs.Item($n)
$xml5 = "<w:p $wNs>" + $listPPr +
    "<w:r>" +
        "<w:rPr><w:rFonts w:hint='default'/><w:lang w:val='pt-BR'/></w:rPr>" +
        "<w:t>draw a strategy first!!!</w:t>" +
    "</w:r>" +
"</w:p>"
$p5.Range.InsertXML($xml5) | Out-Null

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
